$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "Applying functional point light to drawn geometry." (row 31) and
# "Applying functional per pixel spot light to drawn geometry." (row 32)
# as completed on Milestone II, along with rows 30, 33 and 34 which round
# out the Milestone II completion checkmarks for this block.
$ws.Range("E30").Value = "II"
$ws.Range("F30").Value = "X"

$ws.Range("E31").Value = "II"
$ws.Range("F31").Value = "X"

$ws.Range("E32").Value = "II"
$ws.Range("F32").Value = "X"

$ws.Range("E33").Value = "II"
$ws.Range("F33").Value = "X"

$ws.Range("E34").Value = "II"
$ws.Range("F34").Value = "X"

# Update the view to match the latest edit location/selection.
[void]$ws.Range("F33").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
